$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.088304
$ws.Range("H2").Value = 0.264912
$ws.Range("I2").Value = 0.04372337970871547
$ws.Range("J2").Value = 0.04372337970871546
$ws.Range("M2").Value = 34.53319033333333
$ws.Range("N2").Value = 103.599571
$ws.Range("O2").Value = 0.2461870921144496
$ws.Range("P2").Value = 0.2461870921144496
$ws.Range("Q2").Value = 3.049418839194666
$ws.Range("R2").Value = 27.444769552752
$ws.Range("S2").Value = 0.01076413170790459
$ws.Range("T2").Value = 0.01076413170790459

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.088304
$ws.Range("H3").Value = 0.264912
$ws.Range("I3").Value = 0.04372337970871547
$ws.Range("J3").Value = 0.04372337970871546
$ws.Range("O3").Value = 0.575843103803214
$ws.Range("P3").Value = 0.575843103803214
$ws.Range("Q3").Value = 7.132733053045333
$ws.Range("R3").Value = 64.194597477408
$ws.Range("S3").Value = 0.02517780668023318
$ws.Range("T3").Value = 0.02517780668023318

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.088304
$ws.Range("H4").Value = 0.264912
$ws.Range("I4").Value = 0.04372337970871547
$ws.Range("J4").Value = 0.04372337970871546
$ws.Range("O4").Value = 0.1779698040823365
$ws.Range("P4").Value = 0.1779698040823364
$ws.Range("Q4").Value = 2.204439187754667
$ws.Range("R4").Value = 19.839952689792
$ws.Range("S4").Value = 0.007781441320577697
$ws.Range("T4").Value = 0.007781441320577695

# Row 5
$ws.Range("I5").Value = 0.5310748730197871
$ws.Range("J5").Value = 0.531074873019787
$ws.Range("M5").Value = 34.53319033333333
$ws.Range("N5").Value = 103.599571
$ws.Range("O5").Value = 0.2461870921144496
$ws.Range("P5").Value = 0.2461870921144496
$ws.Range("Q5").Value = 37.03898769030066
$ws.Range("R5").Value = 333.350889212706
$ws.Range("S5").Value = 0.1307437786837919
$ws.Range("T5").Value = 0.1307437786837919

# Row 6
$ws.Range("I6").Value = 0.5310748730197871
$ws.Range("J6").Value = 0.531074873019787
$ws.Range("O6").Value = 0.575843103803214
$ws.Range("P6").Value = 0.575843103803214
$ws.Range("S6").Value = 0.3058158032316119
$ws.Range("T6").Value = 0.3058158032316119

# Row 7
$ws.Range("I7").Value = 0.5310748730197871
$ws.Range("J7").Value = 0.531074873019787
$ws.Range("O7").Value = 0.1779698040823365
$ws.Range("P7").Value = 0.1779698040823364
$ws.Range("S7").Value = 0.09451529110438321
$ws.Range("T7").Value = 0.09451529110438318

# Row 8
$ws.Range("I8").Value = 0.4252017472714976
$ws.Range("J8").Value = 0.4252017472714976
$ws.Range("M8").Value = 34.53319033333333
$ws.Range("N8").Value = 103.599571
$ws.Range("O8").Value = 0.2461870921144496
$ws.Range("P8").Value = 0.2461870921144496
$ws.Range("Q8").Value = 29.65503186684667
$ws.Range("R8").Value = 266.89528680162
$ws.Range("S8").Value = 0.1046791817227531
$ws.Range("T8").Value = 0.1046791817227531

# Row 9
$ws.Range("I9").Value = 0.4252017472714976
$ws.Range("J9").Value = 0.4252017472714976
$ws.Range("O9").Value = 0.575843103803214
$ws.Range("P9").Value = 0.575843103803214
$ws.Range("S9").Value = 0.244849493891369
$ws.Range("T9").Value = 0.244849493891369

# Row 10
$ws.Range("I10").Value = 0.4252017472714976
$ws.Range("J10").Value = 0.4252017472714976
$ws.Range("O10").Value = 0.1779698040823365
$ws.Range("P10").Value = 0.1779698040823364
$ws.Range("S10").Value = 0.07567307165737557
$ws.Range("T10").Value = 0.07567307165737555
